$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.946.38"
$ws.Range("E2").Value = "  -1.59%  "
$ws.Range("D3").Value = "3.172.67"
$ws.Range("E3").Value = "  -5.11%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'590.96"
$ws.Range("E5").Value = "  -2.55%  "
$ws.Range("D6").Value = "'134.55"
$ws.Range("E6").Value = "  -6.39%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.170.80"
$ws.Range("E8").Value = "  -5.07%  "
$ws.Range("E9").Value = "  -1.11%  "
$ws.Range("E10").Value = "  -6.65%  "
$ws.Range("D11").Value = "'5.26"
$ws.Range("E11").Value = "  -5.84%  "
$ws.Range("D12").Value = "'0.454"
$ws.Range("E12").Value = "  -3.84%  "
$ws.Range("E13").Value = "  -5.36%  "
$ws.Range("E14").Value = "  -1.06%  "
$ws.Range("D15").Value = "3.691.04"
$ws.Range("E15").Value = "  -5.18%  "
$ws.Range("E16").Value = "  -1.16%  "
$ws.Range("D17").Value = "3.168.00"
$ws.Range("E17").Value = "  -5.22%  "
$ws.Range("D18").Value = "62.869.26"
$ws.Range("E18").Value = "  -1.85%  "
$ws.Range("E19").Value = "  -5.16%  "
$ws.Range("D20").Value = "'461.24"
$ws.Range("E20").Value = "  -4.60%  "
$ws.Range("D21").Value = "'13.91"
$ws.Range("E21").Value = "  -1.99%  "
$ws.Range("D22").Value = "'0.698"
$ws.Range("E22").Value = "  -5.75%  "
$ws.Range("E23").Value = "  -5.07%  "
$ws.Range("D24").Value = "'13.44"
$ws.Range("D25").Value = "'83.00"
$ws.Range("E25").Value = "  -2.49%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("E28").Value = "  -4.63%  "
$ws.Range("E29").Value = "  -5.61%  "
$ws.Range("D30").Value = "'7.72"
$ws.Range("E30").Value = "  -7.28%  "
$ws.Range("E31").Value = "  -6.17%  "
$ws.Range("D32").Value = "'27.14"
$ws.Range("E32").Value = "  -6.38%  "
$ws.Range("E33").Value = "  -4.63%  "
$ws.Range("E34").Value = "  -7.44%  "
$ws.Range("E35").Value = "  -6.59%  "
$ws.Range("E36").Value = "  -4.87%  "
$ws.Range("D37").Value = "'51.31"
$ws.Range("E37").Value = "  -2.37%  "
$ws.Range("D38").Value = "0.0₃0704"
$ws.Range("E38").Value = "  -7.17%  "
$ws.Range("E39").Value = "  -3.63%  "
$ws.Range("D40").Value = "'401.84"
$ws.Range("E40").Value = "  -7.59%  "
$ws.Range("D41").Value = "'8.08"
$ws.Range("E41").Value = "  -3.75%  "
$ws.Range("E42").Value = "  -5.39%  "
$ws.Range("E43").Value = "  -5.98%  "
$ws.Range("D44").Value = "2.793.97"
$ws.Range("E44").Value = "  -11.04%  "
$ws.Range("D45").Value = "'0.250"
$ws.Range("E45").Value = "  -6.74%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("E47").Value = "  -6.94%  "
$ws.Range("D48").Value = "'124.74"
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("D49").Value = "'25.15"
$ws.Range("E49").Value = "  -5.50%  "
$ws.Range("D50").Value = "'34.41"
$ws.Range("E50").Value = "  -6.12%  "
$ws.Range("E51").Value = "  -2.58%  "
